$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# Replace the "Top Hat" outfit (rows 9-12) with a new "Disco" themed outfit:
# Afro / DiscoShirt / DiscoPants / DiscoShoes.
$ws1.Range("C9").Value = "Afro"
$ws1.Range("D9").Value = "Hat"

$ws1.Range("C10").Value = "DiscoShirt"
$ws1.Range("D10").Value = "Shirt"

$ws1.Range("C11").Value = "DiscoPants"
$ws1.Range("D11").Value = "Pants"

$ws1.Range("C12").Value = "DiscoShoes"
$ws1.Range("D12").Value = "Shoes"

# The new, longer item names no longer fit column C - autofit it to the content.
$ws1.Columns.Item(3).AutoFit() | Out-Null

# Move the selection on Sheet1 to reflect where the cursor ended up after the edits.
$ws1.Range("H13").Select() | Out-Null
